$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================
# 1. Header text updates (Volume/Number and report date range)
# =========================================================
$cA8 = $ws.Range("A8")
$chA8 = $cA8.Characters(21, 2)
$chA8.Text = "40"
$chA8.Font.Name = "Andale WT"
$chA8.Font.Size = 10

$cC9 = $ws.Range("C9")
$chC9a = $cC9.Characters(27, 9)
$chC9a.Text = "9/30/2024"
$chC9a.Font.Name = "Andale WT"
$chC9a.Font.Size = 10
$chC9b = $cC9.Characters(47, 9)
$chC9b.Text = "10/6/2024"
$chC9b.Font.Name = "Andale WT"
$chC9b.Font.Size = 10

# =========================================================
# 2. Column E width (now matches columns C/D/F/G/I/J width)
# =========================================================
$ws.Columns.Item(5).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# =========================================================
# 3. Cells changing from "N/A" placeholder text to real numbers
#    (format copied from the untouched Murder row, style 15)
# =========================================================
$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = 2

$ws.Range("I14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 2

$ws.Range("I14").Copy()
$ws.Range("F22").PasteSpecial(-4122)
$ws.Range("F22").Value = 2

$ws.Range("I14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("C27").Value = 3

# =========================================================
# 4. Cells changing from real numbers to "N/A"/"***.* " placeholder text
#    (format + shared text copied from the untouched Murder row, style 14)
# =========================================================
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D20").PasteSpecial(-4122)

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("E14").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "0"
$ws.Range("D14").Copy()
$ws.Range("C31").PasteSpecial(-4122)

# =========================================================
# 5. Plain numeric value updates (style/type unchanged)
# =========================================================
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 14
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = -6.666666666666
$ws.Range("M15").Value = 75
$ws.Range("N15").Value = -17.647058823529
$ws.Range("C16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = -26.666666666666
$ws.Range("I16").Value = 105
$ws.Range("J16").Value = 112
$ws.Range("K16").Value = -6.25
$ws.Range("L16").Value = 17.977528089887
$ws.Range("M16").Value = -23.91304347826
$ws.Range("N16").Value = -82.051282051282
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 5.263157894736
$ws.Range("I17").Value = 202
$ws.Range("J17").Value = 170
$ws.Range("K17").Value = 18.823529411764
$ws.Range("L17").Value = 50.746268656716
$ws.Range("M17").Value = 114.893617021277
$ws.Range("N17").Value = -15.126050420168
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -26.666666666666
$ws.Range("I18").Value = 106
$ws.Range("J18").Value = 168
$ws.Range("K18").Value = -36.904761904761
$ws.Range("L18").Value = -28.378378378378
$ws.Range("M18").Value = -52.466367713004
$ws.Range("N18").Value = -91.246903385631
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -66.666666666666
$ws.Range("F19").Value = 32
$ws.Range("H19").Value = -31.914893617021
$ws.Range("I19").Value = 465
$ws.Range("J19").Value = 512
$ws.Range("K19").Value = -9.1796875
$ws.Range("L19").Value = -21.052631578947
$ws.Range("M19").Value = 48.089171974522
$ws.Range("N19").Value = -20.918367346938
$ws.Range("C20").Value = 5
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 157.142857142857
$ws.Range("I20").Value = 142
$ws.Range("K20").Value = 18.333333333333
$ws.Range("L20").Value = 23.478260869565
$ws.Range("M20").Value = 10.9375
$ws.Range("N20").Value = -89.597069597069
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 25
$ws.Range("E21").Value = -4
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 106
$ws.Range("H21").Value = -10.377358490566
$ws.Range("I21").Value = 1041
$ws.Range("J21").Value = 1098
$ws.Range("K21").Value = -5.191256830601
$ws.Range("L21").Value = -4.931506849315
$ws.Range("M21").Value = 14.270032930845
$ws.Range("N21").Value = -74.052841475573
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 18
$ws.Range("K22").Value = 80
$ws.Range("L22").Value = 80
$ws.Range("M22").Value = -14.285714285714
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 25
$ws.Range("F24").Value = 131
$ws.Range("G24").Value = 99
$ws.Range("H24").Value = 32.323232323232
$ws.Range("I24").Value = 1343
$ws.Range("J24").Value = 1353
$ws.Range("K24").Value = -0.739098300073
$ws.Range("L24").Value = -6.018194541637
$ws.Range("M24").Value = 73.963730569948
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 13
$ws.Range("E25").Value = 23.076923076923
$ws.Range("F25").Value = 70
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = 16.666666666666
$ws.Range("I25").Value = 779
$ws.Range("J25").Value = 783
$ws.Range("K25").Value = -0.510855683269
$ws.Range("L25").Value = -5
$ws.Range("C26").Value = 18
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = -14.285714285714
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 64
$ws.Range("H26").Value = -15.625
$ws.Range("I26").Value = 489
$ws.Range("J26").Value = 465
$ws.Range("K26").Value = 5.16129032258
$ws.Range("L26").Value = 42.151162790697
$ws.Range("M26").Value = 43.823529411764
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 33.333333333333
$ws.Range("I27").Value = 19
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = -5
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 150
$ws.Range("F31").Value = 1

$ws.Range("A1").Select() | Out-Null
